# Update the document's style sheet so the East Asian font used by the
# "Normal" and "Heading" styles changes from "DejaVu Sans" to "Tahoma",
# and so the complex-script ("cs") font is explicitly pinned to
# "DejaVu Sans" on the "List", "Caption" and "Index" styles (mirroring
# what those styles already inherit, but now stated explicitly).

$d = $word.ActiveDocument

# -- East Asian font: DejaVu Sans -> Tahoma -------------------------------
$d.Styles.Item("Normal").Font.NameFarEast  = "Tahoma"
$d.Styles.Item("Heading").Font.NameFarEast = "Tahoma"

# -- Complex-script font: pin explicit "DejaVu Sans" ----------------------
# Font.NameBi maps to <w:rFonts w:cs="..."/> on the style's run properties.
$d.Styles.Item("List").Font.NameBi    = "DejaVu Sans"
$d.Styles.Item("Caption").Font.NameBi = "DejaVu Sans"
$d.Styles.Item("Index").Font.NameBi   = "DejaVu Sans"

Write-Output "styles updated"
